# Apply the "Add files via upload" edit to the active workbook.
# Target sheet is "Sagar Lab" (the active/tab-selected sheet, rId1 -> sheet1.xml).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value edits -------------------------------------------------
# (ordered to match the shared-string append order seen in the target file:
#  Student61, User211, Sakthi, User220)

# B18: student61 -> Student61 (capitalize)
$ws.Range("B18").Value = "Student61"

# B2: user211 -> User211  (capitalize the user name)
$ws.Range("B2").Value = "User211"

# A11: was blank, now holds "Sakthi"
$ws.Range("A11").Value = "Sakthi"

# B11: user220 -> User220 (capitalize)
$ws.Range("B11").Value = "User220"

# A4: "Nageswara Rao" moves out of row 4 (cell cleared, incl. formatting so
# the now-empty cell drops out of the saved sheet entirely)
$ws.Range("A4").Clear()

# A29: "Sakthi S" removed (cell cleared, incl. formatting)
$ws.Range("A29").Clear()

# A42: was blank, now holds "Nageswara Rao" (the name moved down from row 4)
$ws.Range("A42").Value = "Nageswara Rao"

# --- Selection / view state --------------------------------------------
# Select A4 (also clears any scrolled topLeftCell state)
$ws.Range("A4").Select()
